$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "fm40" row (row 3); rows below shift up by one.
$ws.Rows(3).Delete()

# After the shift, rows 5-18 correspond to the sc1-sc14 test cases whose
# "File preparation" status moved from "input files" to "complete".
$ws.Range("F5:F18").Value = "complete"

# sc3-sc8 (now rows 7-12) are now "Currently supported" = yes
$ws.Range("G7:G12").Value = "yes"

# sc5-sc8 (now rows 9-12) are now "Re-implementation" = done
$ws.Range("H9:H12").Value = "done"

# Update the selected range in the sheet view
$ws.Range("B13").Select()
